$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.318.35'
$ws.Cells.Item(2, 5).Value = '  +5.20%  '
$ws.Cells.Item(3, 4).Value = '2.611.76'
$ws.Cells.Item(3, 5).Value = '  +6.65%  '
$ws.Cells.Item(4, 4).Value = '''1.01'
$ws.Cells.Item(4, 5).Value = '  +0.89%  '
$ws.Cells.Item(5, 4).Value = '''506.46'
$ws.Cells.Item(5, 5).Value = '  +2.81%  '
$ws.Cells.Item(6, 4).Value = '''155.62'
$ws.Cells.Item(6, 5).Value = '  +0.56%  '
$ws.Cells.Item(7, 4).Value = '''1.00'
$ws.Cells.Item(7, 5).Value = '  +0.42%  '
$ws.Cells.Item(8, 4).Value = '''0.586'
$ws.Cells.Item(8, 5).Value = '  -4.21%  '
$ws.Cells.Item(9, 4).Value = '2.635.75'
$ws.Cells.Item(9, 5).Value = '  +7.41%  '
$ws.Cells.Item(10, 4).Value = '''6.43'
$ws.Cells.Item(10, 5).Value = '  +4.24%  '
$ws.Cells.Item(11, 4).Value = '''0.104'
$ws.Cells.Item(11, 5).Value = '  +3.70%  '
$ws.Cells.Item(12, 5).Value = '  +2.21%  '
$ws.Cells.Item(13, 5).Value = '  +0.95%  '
$ws.Cells.Item(14, 4).Value = '3.083.77'
$ws.Cells.Item(14, 5).Value = '  +8.37%  '
$ws.Cells.Item(15, 4).Value = '60.414.82'
$ws.Cells.Item(15, 5).Value = '  +5.50%  '
$ws.Cells.Item(16, 4).Value = '''21.69'
$ws.Cells.Item(16, 5).Value = '  +4.43%  '
$ws.Cells.Item(17, 4).Value = '''0.0000139'
$ws.Cells.Item(17, 5).Value = '  +3.93%  '
$ws.Cells.Item(18, 4).Value = '2.627.47'
$ws.Cells.Item(18, 5).Value = '  +7.47%  '
$ws.Cells.Item(19, 4).Value = '''4.77'
$ws.Cells.Item(19, 5).Value = '  +3.21%  '
$ws.Cells.Item(20, 4).Value = '''343.50'
$ws.Cells.Item(20, 5).Value = '  +5.78%  '
$ws.Cells.Item(21, 4).Value = '''10.41'
$ws.Cells.Item(21, 5).Value = '  +3.11%  '
$ws.Cells.Item(22, 4).Value = '''6.16'
$ws.Cells.Item(22, 5).Value = '  +3.42%  '
$ws.Cells.Item(23, 4).Value = '''0.997'
$ws.Cells.Item(23, 5).Value = '  -0.18%  '
$ws.Cells.Item(24, 4).Value = '''5.73'
$ws.Cells.Item(24, 5).Value = '  -0.70%  '
$ws.Cells.Item(25, 4).Value = '''60.14'
$ws.Cells.Item(25, 5).Value = '  +3.93%  '
$ws.Cells.Item(26, 4).Value = '''0.422'
$ws.Cells.Item(26, 5).Value = '  +4.61%  '
$ws.Cells.Item(27, 5).Value = '  +4.36%  '
$ws.Cells.Item(28, 4).Value = '''0.999'
$ws.Cells.Item(28, 5).Value = '  +0.16%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0852'
$ws.Cells.Item(29, 5).Value = '  +6.53%  '
$ws.Cells.Item(30, 4).Value = '''7.51'
$ws.Cells.Item(30, 5).Value = '  +3.02%  '
$ws.Cells.Item(31, 4).Value = '''1.00'
$ws.Cells.Item(31, 5).Value = '  +0.13%  '
$ws.Cells.Item(32, 4).Value = '''156.36'
$ws.Cells.Item(32, 5).Value = '  +4.16%  '
$ws.Cells.Item(33, 4).Value = '''19.35'
$ws.Cells.Item(33, 5).Value = '  +3.06%  '
$ws.Cells.Item(34, 5).Value = '  +2.45%  '
$ws.Cells.Item(35, 4).Value = '''5.71'
$ws.Cells.Item(35, 5).Value = '  +7.42%  '
$ws.Cells.Item(36, 4).Value = '''3.99'
$ws.Cells.Item(36, 5).Value = '  +4.81%  '
$ws.Cells.Item(37, 4).Value = '''1.20'
$ws.Cells.Item(37, 5).Value = '  +6.14%  '
$ws.Cells.Item(38, 2).Value = 'Fetch.AI'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(38, 4).Value = '''0.847'
$ws.Cells.Item(38, 5).Value = '  +2.52%  '
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).Value = '''3.78'
$ws.Cells.Item(39, 5).Value = '  +7.10%  '
$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).Value = '''302.39'
$ws.Cells.Item(40, 5).Value = '  +6.48%  '
$ws.Cells.Item(41, 2).Value = 'SuiNetwork'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(41, 4).Value = '''0.837'
$ws.Cells.Item(41, 5).Value = '  +29.10%  '
$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).Value = '''1.47'
$ws.Cells.Item(42, 5).Value = '  +5.88%  '
$ws.Cells.Item(43, 4).Value = '''35.66'
$ws.Cells.Item(43, 5).Value = '  +4.91%  '
$ws.Cells.Item(44, 4).Value = '''0.628'
$ws.Cells.Item(44, 5).Value = '  +3.88%  '
$ws.Cells.Item(45, 4).Value = '''0.0570'
$ws.Cells.Item(45, 5).Value = '  +6.58%  '
$ws.Cells.Item(46, 4).Value = '''0.100'
$ws.Cells.Item(46, 5).Value = '  -0.39%  '
$ws.Cells.Item(47, 4).Value = '''0.998'
$ws.Cells.Item(47, 5).Value = '  +0.95%  '
$ws.Cells.Item(48, 4).Value = '''19.78'
$ws.Cells.Item(48, 5).Value = '  +10.88%  '
$ws.Cells.Item(49, 4).Value = '''4.86'
$ws.Cells.Item(49, 5).Value = '  +4.92%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).Value = '2.048.78'
$ws.Cells.Item(50, 5).Value = '  +7.47%  '
$ws.Cells.Item(51, 2).Value = 'VeChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(51, 4).Value = '''0.0234'
$ws.Cells.Item(51, 5).Value = '  +2.03%  '
